$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at 59-60; this pushes the existing rows 59-92 down to 61-94.
$ws.Rows("59:60").Insert()

# --- New row 59 ---
$ws.Range("A59").Value = 7
$ws.Range("B59").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C59").Value = "Ñuble"
$ws.Range("D59").Value = 45001
$ws.Range("E59").Value = 16
$ws.Range("F59").Value = "Fruta"
$ws.Range("G59").Value = 100103
$ws.Range("H59").Value = "Frutos de hueso (carozo)"
$ws.Range("I59").Value = 100103002
$ws.Range("J59").Value = "Ciruela"
$ws.Range("K59").Value = "Angeleno"
$ws.Range("L59").Value = "Primera"
$ws.Range("M59").Value = 80
$ws.Range("N59").Value = 11000
$ws.Range("O59").Value = 12000
$ws.Range("P59").Value = 11500
$ws.Range("Q59").Value = '$/bandeja 18 kilos granel'
$ws.Range("R59").Value = "Región de O'Higgins"
$ws.Range("S59").Value = 639
$ws.Range("T59").Value = 18

# --- New row 60 ---
$ws.Range("A60").Value = 7
$ws.Range("B60").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C60").Value = "Ñuble"
$ws.Range("D60").Value = 45001
$ws.Range("E60").Value = 16
$ws.Range("F60").Value = "Fruta"
$ws.Range("G60").Value = 100103
$ws.Range("H60").Value = "Frutos de hueso (carozo)"
$ws.Range("I60").Value = 100103002
$ws.Range("J60").Value = "Ciruela"
$ws.Range("K60").Value = "Angeleno"
$ws.Range("L60").Value = "Segunda"
$ws.Range("M60").Value = 20
$ws.Range("N60").Value = 9500
$ws.Range("O60").Value = 9500
$ws.Range("P60").Value = 9500
$ws.Range("Q60").Value = '$/bandeja 18 kilos granel'
$ws.Range("R60").Value = "Región de O'Higgins"
$ws.Range("S60").Value = 528
$ws.Range("T60").Value = 18
